# fix(mosaic_scripts.xlsx): enable cleanup job
#
# The "active" flag for the "cleanup-mosaic" scheduled job (sheet
# sys_job_ScheduledJob, column E "active", row 2) was stored as the text
# "false". Flip it to the boolean TRUE so the job actually runs, and make
# that sheet the active one (matching the author's saved UI state).

$wb = $excel.ActiveWorkbook

$jobSheet = $wb.Worksheets.Item("sys_job_ScheduledJob")
$jobSheet.Activate()

$jobSheet.Range("E2").Value = $true
$jobSheet.Range("E2").Select() | Out-Null
